$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (member shown in column E, optional "Hoàn thành" status in column F)
$rows = @(
    @{ Row = 18; E = "Bùi Phi Long";           F = "Hoàn thành" },
    @{ Row = 19; E = "Bùi Phi Long";           F = $null },
    @{ Row = 20; E = "Nguyễn Phạm Nhật Minh";  F = $null },
    @{ Row = 21; E = "Nguyễn Phạm Nhật Minh";  F = "Hoàn thành" },
    @{ Row = 22; E = "Nguyễn Phạm Nhật Minh";  F = "Hoàn thành" },
    @{ Row = 23; E = "Nguyễn Phạm Nhật Minh";  F = $null },
    @{ Row = 24; E = "Nguyễn Phạm Nhật Minh";  F = $null },
    @{ Row = 25; E = "Bùi Phi Long";           F = $null },
    @{ Row = 27; E = "Bùi Phi Long";           F = $null },
    @{ Row = 28; E = "Nguyễn Phạm Nhật Minh";  F = $null },
    @{ Row = 29; E = "Nguyễn Phạm Nhật Minh";  F = $null }
)

# Use an already date-formatted cell (D8, short-date style) as the format
# donor so pasting into the blank C/D cells reuses the existing style
# (numFmtId 14) instead of minting a brand-new custom format.
$ws.Cells.Item(8, 4).Copy()

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)
    $ws.Cells.Item($r, 4).PasteSpecial(-4122)
    $ws.Cells.Item($r, 3).Value2 = 45026
    $ws.Cells.Item($r, 4).Value2 = 45240
    $ws.Cells.Item($r, 5).Value = $item.E
    if ($item.F) {
        $ws.Cells.Item($r, 6).Value = $item.F
    }
}

$excel.CutCopyMode = $false
$ws.Range("H13").Select()
